# Refactor data handling and stuff
#
# - Books!D4: "Shaurya" -> "NULL"
# - Active sheet changes from Customers (tab 1) to Books (tab 0)
# - Books sheet selection becomes D8 (and becomes the tab-selected sheet)
# - Customers sheet selection becomes B12 (and loses tab-selected)

$wb = $excel.ActiveWorkbook

$books = $wb.Worksheets.Item("Books")
$customers = $wb.Worksheets.Item("Customers")

# Data edit: Books D4 "Shaurya" -> "NULL"
$books.Range("D4").Value = "NULL"

# Update view/selection state on the no-longer-active sheet first...
$customers.Activate()
$customers.Range("B12").Select()

# ...then activate Books last so it ends up the active/tabSelected sheet.
$books.Activate()
$books.Range("D8").Select()
